# Update the buildings/rooms schedule sheet: reorder/refresh several
# session rows (unit code, classroom number, lecturer, time, delivery mode)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MITS4001"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Jim"

$ws.Range("C3").Value = "MITS4003"
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = "Tom"

$ws.Range("C4").Value = "MITS5503"
$ws.Range("D4").Value = 12
$ws.Range("E4").Value = "Mike"
$ws.Range("F4").Value = "F2F"

$ws.Range("C5").Value = "MITS5507"
$ws.Range("D5").Value = 14
$ws.Range("E5").Value = "Sammy"
$ws.Range("F5").Value = "Online"

$ws.Range("B6").Value = "8:00 AM to 10:00 AM"
$ws.Range("C6").Value = "MITS5501"
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = "Lewis"
$ws.Range("F6").Value = "F2F"

$ws.Range("B7").Value = "8:00 AM to 9:00 AM"
$ws.Range("C7").Value = "MITS5002"
$ws.Range("D7").Value = 7
$ws.Range("E7").Value = "Mitch"
$ws.Range("F7").Value = "Online"

$ws.Range("B8").Value = "9:00 AM to 11:00 AM"
$ws.Range("C8").Value = "MITS6001"
$ws.Range("D8").Value = 16
$ws.Range("E8").Value = "Josh"
$ws.Range("F8").Value = "F2F"

$ws.Range("B9").Value = "9:00 AM to 10:00 AM"
$ws.Range("C9").Value = "MITS4004"
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = "Kat"
$ws.Range("F9").Value = "Online"

$ws.Range("C10").Value = "MITS5502"
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = "Jake"

$ws.Range("C11").Value = "MITS5004"
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = "Sean"

$ws.Range("B15").Value = "2:00 PM to 3:00 PM"
$ws.Range("C15").Value = "MITS6500"
$ws.Range("D15").Value = 21
$ws.Range("E15").Value = "Keno"

$ws.Range("B16").Value = "2:00 PM to 4:00 PM"
$ws.Range("C16").Value = "MITS5003"
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = "Jay"

